# Update sector names in column D to the new taxonomy.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = $ws.UsedRange.Rows.Count }

$map = @{
    "Health Care" = "Healthcare";
    "Information Technology" = "Technology";
    "Consumer Discretionary" = "Consumer Cyclical";
    "Consumer Staples" = "Consumer Defensive";
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
